$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Component Part List: insert a new row (row 29) for a new component -
#    a 24V to 12V DC Converter (1.5A 18W DC DC Converter), and add its
#    hyperlink to a supplier page.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("Component Part List")

$wsComp.Rows.Item(29).Insert()

$wsComp.Range("A29").Value = "Component"
$wsComp.Range("B29").Value = "24V to 12V DC Converter"
$wsComp.Range("C29").Value = "1.5A 18W DC DC Converter"
$wsComp.Range("D29").Value = 1
# E29 is intentionally left blank (price not filled in yet), it keeps the
# numeric style that the row-insert already carried down from row 28.

# Hyperlink on the new component's name to its supplier listing. Adding the
# hyperlink with a TextToDisplay temporarily overwrites the cell text, so
# restore the cell value & its normal "linked part name" style afterwards.
$hl = $wsComp.Hyperlinks.Add($wsComp.Range("B29"), "http://s.click.aliexpress.com/e/bYhNOh2u")
$hl.TextToDisplay = "24V to 12V 1.5A 18W DC Converter"
$wsComp.Range("B29").Value = "24V to 12V DC Converter"
$wsComp.Range("B29").Style = $wsComp.Range("B28").Style

# Widen column C slightly so the longer description still fits.
$wsComp.Columns.Item(3).ColumnWidth = 22

# Resize the small logo picture in the header very slightly.
$logo = $wsComp.Shapes.Item(1)
$logo.LockAspectRatio = $false
$logo.Width = 73.2
$logo.Height = 72

# ---------------------------------------------------------------------------
# 2. V-Slots and Wheels List: update the saved selection to A23 (without
#    making this sheet the active tab).
# ---------------------------------------------------------------------------
$wsSlots = $wb.Worksheets.Item("V-Slots and Wheels List")
$wsSlots.Activate()
$wsSlots.Range("A23").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Make "Component Part List" the active sheet/tab, with B4 selected.
#    (This also clears the "Fasteners List" sheet's previous tabSelected
#    state and moves the active tab index.)
# ---------------------------------------------------------------------------
$wsComp.Activate()
$wsComp.Range("B4").Select() | Out-Null
